# Add two new FAQ entries to the FAQ_MASTER sheet:
#  1. Fix the answer text on the existing "마이크로디그리/CREDIT_INFO" row (F56)
#     so it points at the corrected 소단위전공과정 description.
#  2. Append a brand-new row (row 66) describing how to check grades
#     (성적확인) under the 학사제도 program.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the answer in F56 -----------------------------------------
$ws.Range("F56").Value = "소단위전공과정별로 편성된 과목 3~4과목(8~12학점)을 이수. 원전공 전공과목과 동일한 교과목은 원전공의 이수구분으로, 타 전공의 교과목의 경우 자유선택으로 인정. 교양과목은 원전공의 이수구분을 따름. 소단위전공과정은 이수건수 제한 없음. 이수 시 별도의 이수증을 발급."

# --- 2. Append the new row 66 ----------------------------------------------
$ws.Range("C66").Value = "학사제도"
$ws.Range("D66").Value = "성적확인, 이수학점확인, 내성적확인, 수강신청확인, 내수강신청확인"
$ws.Range("F66").Value = "현재까지의 성적, 이수학점 등 확인은 학사시스템 https://info.hknu.ac.kr 또는 학사지원팀 031-670-5032 연락"

# Match styling used by the surrounding rows (style index "2": vertical
# center + shrink-to-fit, same as the rest of the table). Order matters:
# ShrinkToFit first, then VerticalAlignment, so the engine resolves to the
# already-existing combined style instead of minting a new one.
foreach ($addr in "C66", "D66", "F66") {
    $rng = $ws.Range($addr)
    $rng.ShrinkToFit = $true
    $rng.VerticalAlignment = -4108
}

# --- 3. Update the view so the newly-added row is visible/selected --------
$ws.Range("F66").Select()
$excel.ActiveWindow.ScrollRow = 51
